# Append 7 new "Noun" log rows (rows 10-16) to the sheet, mirroring the
# existing rows' layout. This represents new entries written by the
# (now thread-safe) log textbox.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(42601.976909722223, "Noun", 3125, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42601.988217592596, "Noun", 3111, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42601.98877314815,  "Noun", 3027, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42601.990972222222, "Noun", 3200, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42601.994826388887, "Noun", 3117, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42601.997800925928, "Noun", 2794, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(42601.99895833333,  "Noun", 3040, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42)
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
